# Scheduled runner update: refresh scraped market-board price/profit
# figures across the Faerie_Profits leve-profit sheets (one per crafting
# job). Each block below targets a single row's currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ /
# LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns (H:N) with freshly
# pulled values. A few rows flip which of LeveProfitNQ/HQ (M/N) is
# populated, so those cells are cleared rather than merely overwritten.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3857.5833
$ws.Range("I113").Value = 2723
$ws.Range("J113").Value = 4424.875
$ws.Range("K113").Value = 2723
$ws.Range("L113").Value = 4424.875
$ws.Range("M113").Value = 531
$ws.Range("N113").Value = -10932.875

$ws.Range("H116").Value = 2418.182
$ws.Range("I116").Value = 2410
$ws.Range("K116").Value = 2410
$ws.Range("M116").Value = 1032

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11253.281
$ws.Range("I32").Value = 12147.679
$ws.Range("J32").Value = 4992.5
$ws.Range("K32").Value = 12147.679
$ws.Range("L32").Value = 4992.5
$ws.Range("M32").Value = -11860.679
$ws.Range("N32").Value = -5566.5

$ws.Range("H88").Value = 2246.7778
$ws.Range("J88").Value = 2031.5714
$ws.Range("L88").Value = 2031.5714
$ws.Range("N88").Value = -2843.5714

$ws.Range("H91").Value = 2246.7778
$ws.Range("J91").Value = 2031.5714
$ws.Range("L91").Value = 2031.5714
$ws.Range("N91").Value = -4839.5714

$ws.Range("H97").Value = 884.875
$ws.Range("I97").Value = 898.4286
$ws.Range("J97").Value = 790
$ws.Range("K97").Value = 898.4286
$ws.Range("L97").Value = 790
$ws.Range("M97").Value = -402.4286
$ws.Range("N97").Value = -1782

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2797.682
$ws.Range("I20").Value = 2574.1538
$ws.Range("J20").Value = 3120.5557
$ws.Range("K20").Value = 2574.1538
$ws.Range("L20").Value = 3120.5557
$ws.Range("M20").Value = -2327.1538
$ws.Range("N20").Value = -3614.5557

$ws.Range("H29").Value = 13098.5
$ws.Range("I29").Value = 13098.5
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 13098.5
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -12809.5
$ws.Range("N29").ClearContents()

$ws.Range("H86").Value = 3511112.5
$ws.Range("J86").Value = 1647.5
$ws.Range("L86").Value = 1647.5
$ws.Range("N86").Value = -3893.5

$ws.Range("H89").Value = 3511112.5
$ws.Range("J89").Value = 1647.5
$ws.Range("L89").Value = 8237.5
$ws.Range("N89").Value = -19469.5

$ws.Range("H134").Value = 7183.6807
$ws.Range("I134").Value = 2602.1875
$ws.Range("J134").Value = 9548.322
$ws.Range("K134").Value = 7806.5625
$ws.Range("L134").Value = 28644.966
$ws.Range("M134").Value = -5271.5625
$ws.Range("N134").Value = -33714.966

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2211.8462
$ws.Range("I58").Value = 2466.5715
$ws.Range("J58").Value = 1914.6666
$ws.Range("K58").Value = 2466.5715
$ws.Range("L58").Value = 1914.6666
$ws.Range("M58").Value = -2263.5715
$ws.Range("N58").Value = -2320.6666

$ws.Range("H86").Value = 5111.316
$ws.Range("I86").Value = 3727
$ws.Range("K86").Value = 3727
$ws.Range("M86").Value = -2604

$ws.Range("H89").Value = 5111.316
$ws.Range("I89").Value = 3727
$ws.Range("K89").Value = 18635
$ws.Range("M89").Value = -13019

$ws.Range("H136").Value = 2211.8462
$ws.Range("I136").Value = 2466.5715
$ws.Range("J136").Value = 1914.6666
$ws.Range("K136").Value = 7399.7145
$ws.Range("L136").Value = 5743.9998
$ws.Range("M136").Value = -4849.7145
$ws.Range("N136").Value = -10843.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 9293553
$ws.Range("I4").Value = 19952692
$ws.Range("J4").Value = 3710195
$ws.Range("K4").Value = 59858076
$ws.Range("L4").Value = 11130585
$ws.Range("M4").Value = -59857964
$ws.Range("N4").Value = -11130809

$ws.Range("H107").Value = 4388.923
$ws.Range("I107").Value = 396.5
$ws.Range("J107").Value = 5114.8184
$ws.Range("K107").Value = 1189.5
$ws.Range("L107").Value = 15344.4552
$ws.Range("M107").Value = 730.5
$ws.Range("N107").Value = -19184.4552

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 39598
$ws.Range("I32").Value = 39399.5
$ws.Range("J32").Value = 39995
$ws.Range("K32").Value = 39399.5
$ws.Range("L32").Value = 39995
$ws.Range("M32").Value = -39103.5
$ws.Range("N32").Value = -40587

$ws.Range("H34").Value = 189914.25
$ws.Range("J34").Value = 186833.33
$ws.Range("L34").Value = 186833.33
$ws.Range("N34").Value = -187369.33

$ws.Range("H42").Value = 99900
$ws.Range("I42").Value = 99900
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 99900
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -99415
$ws.Range("N42").ClearContents()

$ws.Range("H70").Value = 8183.3335
$ws.Range("I70").Value = 7775
$ws.Range("K70").Value = 7775
$ws.Range("M70").Value = -7505

$ws.Range("H73").Value = 8183.3335
$ws.Range("I73").Value = 7775
$ws.Range("K73").Value = 7775
$ws.Range("M73").Value = -6839

$ws.Range("H76").Value = 189914.25
$ws.Range("J76").Value = 186833.33
$ws.Range("L76").Value = 186833.33
$ws.Range("N76").Value = -187463.33

$ws.Range("H79").Value = 189914.25
$ws.Range("J79").Value = 186833.33
$ws.Range("L79").Value = 186833.33
$ws.Range("N79").Value = -189017.33

$ws.Range("H115").Value = 99900
$ws.Range("I115").Value = 99900
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 99900
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -98725
$ws.Range("N115").ClearContents()

$ws.Range("H122").Value = 1851.75
$ws.Range("I122").Value = 1447.15
$ws.Range("J122").Value = 3874.75
$ws.Range("K122").Value = 4341.450000000001
$ws.Range("L122").Value = 11624.25
$ws.Range("M122").Value = -1891.450000000001
$ws.Range("N122").Value = -16524.25

$ws.Range("H125").Value = 20000
$ws.Range("J125").Value = 20000
$ws.Range("L125").Value = 20000
$ws.Range("N125").Value = -24920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2570.111
$ws.Range("J22").Value = 2798.25
$ws.Range("L22").Value = 2798.25
$ws.Range("N22").Value = -3388.25

$ws.Range("H27").Value = 2570.111
$ws.Range("J27").Value = 2798.25
$ws.Range("L27").Value = 2798.25
$ws.Range("N27").Value = -3012.25

$ws.Range("H40").Value = 4648.6665
$ws.Range("I40").Value = 4653.0645
$ws.Range("K40").Value = 4653.0645
$ws.Range("M40").Value = -4517.0645

$ws.Range("H55").Value = 262.82352
$ws.Range("I55").Value = 277
$ws.Range("J55").Value = 242.57143
$ws.Range("K55").Value = 277
$ws.Range("L55").Value = 242.57143
$ws.Range("M55").Value = -104
$ws.Range("N55").Value = -588.57143

$ws.Range("H136").Value = 4056.4
$ws.Range("I136").Value = 3390
$ws.Range("K136").Value = 10170
$ws.Range("M136").Value = -7620

$ws.Range("H141").Value = 65000
$ws.Range("I141").Value = 45000
$ws.Range("J141").Value = 75000
$ws.Range("K141").Value = 45000
$ws.Range("L141").Value = 75000
$ws.Range("M141").Value = -39820
$ws.Range("N141").Value = -85360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 11465.211
$ws.Range("I136").Value = 13302.125
$ws.Range("J136").Value = 1668.3334
$ws.Range("K136").Value = 39906.375
$ws.Range("L136").Value = 5005.0002
$ws.Range("M136").Value = -37356.375
$ws.Range("N136").Value = -10105.0002

